$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("LiveData")

$ws.Range("C3").Value = 160747
$ws.Range("C4").Value = 151771
$ws.Range("C7").Value = 5.58
$ws.Range("C8").Value = 64.36
